$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.071.44"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  +0.17%  "
$ws.Cells.Item(2, 5).ClearFormats()

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.911.62"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  +0.52%  "
$ws.Cells.Item(3, 5).ClearFormats()

$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$ws.Cells.Item(4, 5).ClearFormats()

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.7867"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  +6.13%  "
$ws.Cells.Item(5, 5).ClearFormats()

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "243.13"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  -0.17%  "
$ws.Cells.Item(6, 5).ClearFormats()

$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  +0.04%  "
$ws.Cells.Item(7, 5).ClearFormats()

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3168"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  +3.30%  "
$ws.Cells.Item(8, 5).ClearFormats()

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "26.24"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  +0.46%  "
$ws.Cells.Item(9, 5).ClearFormats()

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.06927"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  +0.38%  "
$ws.Cells.Item(10, 5).ClearFormats()

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07983"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  -0.56%  "
$ws.Cells.Item(11, 5).ClearFormats()

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.914.46"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  +0.66%  "
$ws.Cells.Item(12, 5).ClearFormats()

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.7461"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  -2.21%  "
$ws.Cells.Item(13, 5).ClearFormats()

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.223"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  -0.18%  "
$ws.Cells.Item(14, 5).ClearFormats()

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "93.37"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  +2.24%  "
$ws.Cells.Item(15, 5).ClearFormats()

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "30.095.78"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "  +0.23%  "
$ws.Cells.Item(16, 5).ClearFormats()

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.01"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -0.10%  "
$ws.Cells.Item(17, 5).ClearFormats()

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "5.917"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  -2.88%  "
$ws.Cells.Item(18, 5).ClearFormats()

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "247.10"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  +4.00%  "
$ws.Cells.Item(19, 5).ClearFormats()

$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  +0.32%  "
$ws.Cells.Item(20, 5).ClearFormats()

$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  +0.04%  "
$ws.Cells.Item(21, 5).ClearFormats()

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.001"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  +0.14%  "
$ws.Cells.Item(22, 5).ClearFormats()

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.912"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  -2.01%  "
$ws.Cells.Item(23, 5).ClearFormats()

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "169.39"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  +1.82%  "
$ws.Cells.Item(24, 5).ClearFormats()

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.310"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  -0.06%  "
$ws.Cells.Item(25, 5).ClearFormats()

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1380"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  +9.43%  "
$ws.Cells.Item(26, 5).ClearFormats()

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.92"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  +0.51%  "
$ws.Cells.Item(27, 5).ClearFormats()

$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  +0.47%  "
$ws.Cells.Item(28, 5).ClearFormats()

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.375"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  +1.77%  "
$ws.Cells.Item(29, 5).ClearFormats()

$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  -0.82%  "
$ws.Cells.Item(30, 5).ClearFormats()

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.340"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  +1.19%  "
$ws.Cells.Item(31, 5).ClearFormats()

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.05745"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  +8.01%  "
$ws.Cells.Item(32, 5).ClearFormats()

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.117"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  +1.95%  "
$ws.Cells.Item(33, 5).ClearFormats()

$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  -2.43%  "
$ws.Cells.Item(34, 5).ClearFormats()

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7367"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  -0.01%  "
$ws.Cells.Item(35, 5).ClearFormats()

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.723"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  +0.14%  "
$ws.Cells.Item(36, 5).ClearFormats()

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01923"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  -1.20%  "
$ws.Cells.Item(37, 5).ClearFormats()

$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  +0.05%  "
$ws.Cells.Item(38, 5).ClearFormats()

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "6.164"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  -1.74%  "
$ws.Cells.Item(39, 5).ClearFormats()

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.4449"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  -0.07%  "
$ws.Cells.Item(40, 5).ClearFormats()

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "72.64"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  -0.64%  "
$ws.Cells.Item(41, 5).ClearFormats()

$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  +0.03%  "
$ws.Cells.Item(42, 5).ClearFormats()

$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  -2.99%  "
$ws.Cells.Item(43, 5).ClearFormats()

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.8356"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  +0.24%  "
$ws.Cells.Item(44, 5).ClearFormats()

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "7.567"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  -0.77%  "
$ws.Cells.Item(45, 5).ClearFormats()

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "100.57"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  -0.77%  "
$ws.Cells.Item(46, 5).ClearFormats()

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.817"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  +0.12%  "
$ws.Cells.Item(47, 5).ClearFormats()

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "992.82"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  +8.55%  "
$ws.Cells.Item(48, 5).ClearFormats()

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.056.56"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "  -0.05%  "
$ws.Cells.Item(49, 5).ClearFormats()

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "36.19"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  -0.86%  "
$ws.Cells.Item(50, 5).ClearFormats()

$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  +2.84%  "
$ws.Cells.Item(51, 5).ClearFormats()
